## Generate Report for Handoff
## Refreshes the latest-handoff row on each sheet and appends the
## two newly discovered handoff rows (Overview + per-locale detail).
$wb = $excel.ActiveWorkbook

# --- Sheet1: Overview ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("A2").Hyperlinks.Delete()   # existing link refreshed in place

# Row 2
$ws1.Range("A2").Value = "0022afde-b463-4dfb-9a83-c3df912330d7.png"
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"
$ws1.Range("D2").Value = "2016-47-19 18:47:57"
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c4ddd76785156a6a3cb128ffbdc2856dbd15de5e/e2e/0022afde-b463-4dfb-9a83-c3df912330d7.png")

# Row 3
$ws1.Range("A3").Value = "552946cb-d4b4-47e2-b00e-a92d320443d6.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Range("D3").Value = "2016-47-19 18:47:57"
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c4ddd76785156a6a3cb128ffbdc2856dbd15de5e/e2e/552946cb-d4b4-47e2-b00e-a92d320443d6.md")

# Row 4
$ws1.Range("A4").Value = "76c76c9d-34db-48a9-906d-9cc3106487fd.png"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"
$ws1.Range("D4").Value = "2016-47-19 18:47:57"
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/c4ddd76785156a6a3cb128ffbdc2856dbd15de5e/e2e/76c76c9d-34db-48a9-906d-9cc3106487fd.png")

# --- Sheet2: zh-cn ---
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A2").Hyperlinks.Delete()   # existing links refreshed in place
$ws2.Range("B2").Hyperlinks.Delete()
$ws2.Range("D2").Hyperlinks.Delete()

# Row 2
$ws2.Range("A2").Value = "0022afde-b463-4dfb-9a83-c3df912330d7.png"
$ws2.Range("B2").Value = ".png"
$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("D2").Value = "81b0ef128c925bb4e3f24a3903a36bb772af7557.png"
$ws2.Range("E2").Value = "2016-03-19 18:47:55"
$ws2.Range("H2").Value = "0001-01-01 00:00:00"
$ws2.Range("I2").Value = "IsDependency"
$ws2.Range("J2").Value = "e2e\552946cb-d4b4-47e2-b00e-a92d320443d6.md"
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c4ddd76785156a6a3cb128ffbdc2856dbd15de5e/e2e/0022afde-b463-4dfb-9a83-c3df912330d7.png")
$ws2.Hyperlinks.Add($ws2.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/c4ddd76785156a6a3cb128ffbdc2856dbd15de5e/e2e/0022afde-b463-4dfb-9a83-c3df912330d7.png")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f80ecc335243e85aaf2a2640b9080304aed4652c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/81b0ef128c925bb4e3f24a3903a36bb772af7557.png")

# Row 3
$ws2.Range("A3").Value = "552946cb-d4b4-47e2-b00e-a92d320443d6.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "552946cb-d4b4-47e2-b00e-a92d320443d6.11af609bae015369faf8f2dcb1b56facecc154cf.zh-cn.xlf"
$ws2.Range("E3").Value = "2016-03-19 18:47:55"
$ws2.Range("H3").Value = "0001-01-01 00:00:00"
$ws2.Range("I3").Value = "Include"
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c4ddd76785156a6a3cb128ffbdc2856dbd15de5e/e2e/552946cb-d4b4-47e2-b00e-a92d320443d6.md")
$ws2.Hyperlinks.Add($ws2.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/c4ddd76785156a6a3cb128ffbdc2856dbd15de5e/e2e/552946cb-d4b4-47e2-b00e-a92d320443d6.md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f80ecc335243e85aaf2a2640b9080304aed4652c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/552946cb-d4b4-47e2-b00e-a92d320443d6.11af609bae015369faf8f2dcb1b56facecc154cf.zh-cn.xlf")

# Row 4
$ws2.Range("A4").Value = "76c76c9d-34db-48a9-906d-9cc3106487fd.png"
$ws2.Range("B4").Value = ".png"
$ws2.Range("C4").Value = "Ready for handoff"
$ws2.Range("D4").Value = "a0c2c2066b6bcc60255f3708a44552e645af080c.png"
$ws2.Range("E4").Value = "2016-03-19 18:47:55"
$ws2.Range("H4").Value = "0001-01-01 00:00:00"
$ws2.Range("I4").Value = "IsDependency"
$ws2.Range("J4").Value = "e2e\552946cb-d4b4-47e2-b00e-a92d320443d6.md"
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/c4ddd76785156a6a3cb128ffbdc2856dbd15de5e/e2e/76c76c9d-34db-48a9-906d-9cc3106487fd.png")
$ws2.Hyperlinks.Add($ws2.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/c4ddd76785156a6a3cb128ffbdc2856dbd15de5e/e2e/76c76c9d-34db-48a9-906d-9cc3106487fd.png")
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f80ecc335243e85aaf2a2640b9080304aed4652c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a0c2c2066b6bcc60255f3708a44552e645af080c.png")

# --- Sheet3: de-de ---
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A2").Hyperlinks.Delete()   # existing links refreshed in place
$ws3.Range("B2").Hyperlinks.Delete()
$ws3.Range("D2").Hyperlinks.Delete()

# Row 2
$ws3.Range("A2").Value = "0022afde-b463-4dfb-9a83-c3df912330d7.png"
$ws3.Range("B2").Value = ".png"
$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("D2").Value = "81b0ef128c925bb4e3f24a3903a36bb772af7557.png"
$ws3.Range("E2").Value = "2016-03-19 18:47:57"
$ws3.Range("H2").Value = "0001-01-01 00:00:00"
$ws3.Range("I2").Value = "IsDependency"
$ws3.Range("J2").Value = "e2e\552946cb-d4b4-47e2-b00e-a92d320443d6.md"
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c4ddd76785156a6a3cb128ffbdc2856dbd15de5e/e2e/0022afde-b463-4dfb-9a83-c3df912330d7.png")
$ws3.Hyperlinks.Add($ws3.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/c4ddd76785156a6a3cb128ffbdc2856dbd15de5e/e2e/0022afde-b463-4dfb-9a83-c3df912330d7.png")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/def3881369abbf632a5c330fe4624b940ab43bf6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/81b0ef128c925bb4e3f24a3903a36bb772af7557.png")

# Row 3
$ws3.Range("A3").Value = "552946cb-d4b4-47e2-b00e-a92d320443d6.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "552946cb-d4b4-47e2-b00e-a92d320443d6.11af609bae015369faf8f2dcb1b56facecc154cf.de-de.xlf"
$ws3.Range("E3").Value = "2016-03-19 18:47:57"
$ws3.Range("H3").Value = "0001-01-01 00:00:00"
$ws3.Range("I3").Value = "Include"
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c4ddd76785156a6a3cb128ffbdc2856dbd15de5e/e2e/552946cb-d4b4-47e2-b00e-a92d320443d6.md")
$ws3.Hyperlinks.Add($ws3.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/c4ddd76785156a6a3cb128ffbdc2856dbd15de5e/e2e/552946cb-d4b4-47e2-b00e-a92d320443d6.md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/def3881369abbf632a5c330fe4624b940ab43bf6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/552946cb-d4b4-47e2-b00e-a92d320443d6.11af609bae015369faf8f2dcb1b56facecc154cf.de-de.xlf")

# Row 4
$ws3.Range("A4").Value = "76c76c9d-34db-48a9-906d-9cc3106487fd.png"
$ws3.Range("B4").Value = ".png"
$ws3.Range("C4").Value = "Ready for handoff"
$ws3.Range("D4").Value = "a0c2c2066b6bcc60255f3708a44552e645af080c.png"
$ws3.Range("E4").Value = "2016-03-19 18:47:57"
$ws3.Range("H4").Value = "0001-01-01 00:00:00"
$ws3.Range("I4").Value = "IsDependency"
$ws3.Range("J4").Value = "e2e\552946cb-d4b4-47e2-b00e-a92d320443d6.md"
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/c4ddd76785156a6a3cb128ffbdc2856dbd15de5e/e2e/76c76c9d-34db-48a9-906d-9cc3106487fd.png")
$ws3.Hyperlinks.Add($ws3.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/c4ddd76785156a6a3cb128ffbdc2856dbd15de5e/e2e/76c76c9d-34db-48a9-906d-9cc3106487fd.png")
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/def3881369abbf632a5c330fe4624b940ab43bf6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a0c2c2066b6bcc60255f3708a44552e645af080c.png")
